# "Generate Report for Handback"
#
# The localization round-trip for both e2e docs (c83a7963... and
# da5e691d...) has come back "in sync with en-US", so:
#   - the Overview sheet's per-language status columns flip from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - each language sheet (zh-cn / de-de) gets its "Latest Target File" /
#     "Latest Handback File" / "Latest Handback DateTime" columns filled
#     in (they were blank placeholders before the first handback), with
#     the Latest Target File cell turned into a clickable link back to
#     the source doc, same as column A
#   - a couple of columns get widened so the newly-populated long
#     filenames aren't truncated

$wb = $excel.ActiveWorkbook

$c83aMd  = "c83a7963-49ac-44ce-ba7d-43f6ffdb3b24.md"
$c83aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67cb649e67b0d62b4709669f6aa2183e3dd3f34b/e2e/c83a7963-49ac-44ce-ba7d-43f6ffdb3b24.md"
$da5eMd  = "da5e691d-f6b3-4ecb-a132-7b63c4c63144.md"
$da5eUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67cb649e67b0d62b4709669f6aa2183e3dd3f34b/e2e/da5e691d-f6b3-4ecb-a132-7b63c4c63144.md"

# ---------------------------------------------------------------------
# Overview sheet: status text for both languages, both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn / de-de status columns widened to fit the longer status text
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $c83aUrl, "", "", $c83aMd) | Out-Null
$zhcn.Range("J2").Value = "c83a7963-49ac-44ce-ba7d-43f6ffdb3b24.088220d5f050c1c756df7e350955dbd304a607c9.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-06 08:37:50"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $da5eUrl, "", "", $da5eMd) | Out-Null
$zhcn.Range("J3").Value = "da5e691d-f6b3-4ecb-a132-7b63c4c63144.98d586997bd5e18045493745a2ae8d971c1d8f7f.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-06 08:37:50"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

$dede.Hyperlinks.Add($dede.Range("I2"), $c83aUrl, "", "", $c83aMd) | Out-Null
$dede.Range("J2").Value = "c83a7963-49ac-44ce-ba7d-43f6ffdb3b24.088220d5f050c1c756df7e350955dbd304a607c9.de-de.xlf"
$dede.Range("K2").Value = "2016-09-06 08:38:19"

$dede.Hyperlinks.Add($dede.Range("I3"), $da5eUrl, "", "", $da5eMd) | Out-Null
$dede.Range("J3").Value = "da5e691d-f6b3-4ecb-a132-7b63c4c63144.98d586997bd5e18045493745a2ae8d971c1d8f7f.de-de.xlf"
$dede.Range("K3").Value = "2016-09-06 08:38:19"
